# se arreglo item que daba error con el Span en Hogar
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DatosCuenta")

# Update the Span item text from "SmokeSiete" to "SmokeOcho"
$ws.Range("A2").Value = "SmokeOcho"
$ws.Range("B2").Value = "SmokeOcho"

# Update the related numeric values
$ws.Range("C2").Value = 21546999
$ws.Range("D2").Value = 141

# Move the active selection to D6
$ws.Range("D6").Select()
